$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.659.36'
$ws.Range("E2").Value = '  +2.25%  '

$ws.Range("D3").Value = '1.873.20'
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.78'
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4669'
$ws.Range("E7").Value = '  +1.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3894'
$ws.Range("E8").Value = '  +0.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07889'
$ws.Range("E9").Value = '  +0.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9742'
$ws.Range("E10").Value = '  +1.60%  '

$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("D12").Value = '1.888.62'
$ws.Range("E12").Value = '  +3.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.997'
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.713'
$ws.Range("E14").Value = '  +1.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06981'
$ws.Range("E15").Value = '  +3.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.25'
$ws.Range("E16").Value = '  +1.81%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.57%  '

$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.86'
$ws.Range("E19").Value = '  +1.27%  '

$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("D21").Value = '28.662.44'
$ws.Range("E21").Value = '  +2.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.305'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.03'
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.115'
$ws.Range("E24").Value = '  +0.60%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.144.67'
$ws.Range("E25").Value = '  +3.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.51'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.27'
$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.752'
$ws.Range("E28").Value = '  +0.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.988'
$ws.Range("E29").Value = '  +0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.42'
$ws.Range("E30").Value = '  +2.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09364'
$ws.Range("E31").Value = '  +1.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9188'
$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.277'
$ws.Range("E33").Value = '  -0.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.341'
$ws.Range("E34").Value = '  +1.91%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.348'
$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05806'
$ws.Range("E36").Value = '  -1.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02100'
$ws.Range("E37").Value = '  -2.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.147'
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.772'
$ws.Range("E39").Value = '  +0.60%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5624'
$ws.Range("E40").Value = '  +0.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1786'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.790'
$ws.Range("E42").Value = '  -1.13%  '

$ws.Range("E43").Value = '  +2.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.75'
$ws.Range("E44").Value = '  +1.63%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5318'
$ws.Range("E45").Value = '  +0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.163'
$ws.Range("E46").Value = '  -4.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.827'
$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.41'
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.064'
$ws.Range("E49").Value = '  -3.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.366'
$ws.Range("E50").Value = '  +2.06%  '

$ws.Range("E51").Value = '  +0.51%  '
